$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Gold"
$ws.Range("B2").Value = 50000

# "duration" for the Gold package is the text "1" (not the number 1), so
# route the literal through a TEXT() formula and then flatten it down to a
# plain value via copy / paste-special so the cell ends up holding a real
# text value (no residual formula, no quote-prefix formatting).
$ws.Range("C2").Formula = '=TEXT(1,"0")'
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4163) | Out-Null

$ws.Range("D2").Value = "Group Training"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "6 days/week"
$ws.Range("G2").Value = "6 days/week"
$ws.Range("H2").Value = "All Day Access"
